$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "27.154.25"
$ws.Range("E2").Value = "  +0.40%  "
$ws.Range("D3").Value = "1.683.86"
$ws.Range("E3").Value = "  +0.02%  "
$ws.Range("E4").Value = "  +0.19%  "
$ws.Range("D5").Value = "215.87"
$ws.Range("E5").Value = "  -0.16%  "
$ws.Range("D6").Value = "0.518"
$ws.Range("E6").Value = "  +0.05%  "
$ws.Range("E7").Value = "  +0.11%  "
$ws.Range("D8").Value = "22.87"
$ws.Range("E8").Value = "  +5.96%  "
$ws.Range("E9").Value = "  +2.45%  "
$ws.Range("E10").Value = "  +0.36%  "
$ws.Range("D11").Value = "0.0890"
$ws.Range("E11").Value = "  -0.10%  "
$ws.Range("D12").Value = "1.922.40"
$ws.Range("E12").Value = "  +0.08%  "
$ws.Range("D13").Value = "1.692.10"
$ws.Range("E13").Value = "  +0.56%  "
$ws.Range("E14").Value = "  +2.00%  "
$ws.Range("D15").Value = "0.559"
$ws.Range("E15").Value = "  +4.53%  "
$ws.Range("D16").Value = "66.86"
$ws.Range("E16").Value = "  +0.70%  "
$ws.Range("D17").Value = "27.167.89"
$ws.Range("E17").Value = "  +0.43%  "
$ws.Range("D18").Value = "235.65"
$ws.Range("E18").Value = "  -0.52%  "
$ws.Range("D19").Value = "7.96"
$ws.Range("E19").Value = "  -2.85%  "
$ws.Range("E20").Value = "  +0.60%  "
$ws.Range("E21").Value = "  +0.12%  "
$ws.Range("E22").Value = "  +1.94%  "
$ws.Range("D23").Value = "9.56"
$ws.Range("E23").Value = "  +2.94%  "
$ws.Range("E24").Value = "  -2.06%  "
$ws.Range("D25").Value = "146.87"
$ws.Range("E25").Value = "  -0.11%  "
$ws.Range("E26").Value = "  +1.60%  "
$ws.Range("D27").Value = "16.41"
$ws.Range("E27").Value = "  -2.34%  "
$ws.Range("E28").Value = "  +0.05%  "
$ws.Range("E29").Value = "  -0.02%  "
$ws.Range("E30").Value = "  +0.90%  "
$ws.Range("E31").Value = "  +0.03%  "
$ws.Range("E32").Value = "  +0.64%  "
$ws.Range("D33").Value = "1.546.12"
$ws.Range("E33").Value = "  +1.36%  "
$ws.Range("E34").Value = "  +1.66%  "
$ws.Range("E35").Value = "  -2.07%  "
$ws.Range("D36").Value = "0.603"
$ws.Range("E36").Value = "  +2.06%  "
$ws.Range("E37").Value = "  +2.57%  "
$ws.Range("E38").Value = "  -0.47%  "
$ws.Range("E39").Value = "  -1.28%  "
$ws.Range("E40").Value = "  +2.47%  "
$ws.Range("D41").Value = "5.76"
$ws.Range("E41").Value = "  +0.00%  "
$ws.Range("D42").Value = "69.04"
$ws.Range("E42").Value = "  +1.33%  "
$ws.Range("E44").Value = "  -0.45%  "
$ws.Range("D45").Value = "1.829.03"
$ws.Range("E45").Value = "  +0.28%  "
$ws.Range("D46").Value = "0.792"
$ws.Range("E46").Value = "  +1.28%  "
$ws.Range("D47").Value = "90.07"
$ws.Range("E47").Value = "  -0.24%  "
$ws.Range("E48").Value = "  +6.95%  "
$ws.Range("E49").Value = "  +5.54%  "
$ws.Range("D50").Value = "8.25"
$ws.Range("E50").Value = "  +4.13%  "
$ws.Range("E51").Value = "  -0.92%  "
